# Updated INS datasets page function
#
# The "TabQuery" SQL text stored for the Diagnosis tab (row 4 / cell B4)
# and the Survival tab (row 7 / cell B7) had their ORDER BY clause changed
# so that the participant id is cast to text before sorting, i.e.:
#
#   ORDER BY
#       prt.participant_id ASC
#
# becomes
#
#   ORDER BY
#       CAST(prt.participant_id AS TEXT) ASC
#
# Everything else about those two SQL statements (and the rest of the
# sheet) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldOrderBy = "ORDER BY `n    prt.participant_id ASC"
$newOrderBy = "ORDER BY `n    CAST(prt.participant_id AS TEXT) ASC"

# DiagnosisTab row -> B4 holds the Diagnosis query text
$ws.Range("B4").Replace($oldOrderBy, $newOrderBy)

# SurvivalTab row -> B7 holds the Survival query text
$ws.Range("B7").Replace($oldOrderBy, $newOrderBy)

# The workbook was left scrolled down a couple of rows with C7 as the
# active/selected cell (previously the view was parked at C5).
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select()
